$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before column J (10) to make room for the new
#    "RequestContent" column. This shifts the old J/K/L (SchemaValidation,
#    Csvson, Tags) columns one place to the right, to K/L/M, carrying their
#    styles along - exactly mirroring Excel's native "Insert Column" command.
# ---------------------------------------------------------------------------
$ws.Columns.Item(10).Insert()

# ---------------------------------------------------------------------------
# 2. Row 3 - "Test one-POST"
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Test one-POST"
$ws.Range("B3").Value = "REST"
$ws.Range("C3").Value = "ep"
$ws.Range("D3").Value = "Post users"
$ws.Range("E3").Value = "https://reqres.in/api/users"
$ws.Range("F3").Value = "POST"
$ws.Range("G3").Value = "application/json"
$ws.Range("H3").Value = 201
$ws.Range("I3").Value = "name=venom"

# New column header (introduced at this point to match the workbook's
# shared-string ordering)
$ws.Range("J1").Value = "RequestContent"

$ws.Range("J3").Value = "{`n    ""name"": ""venom"",`n    ""job"": ""snake""`n}"
$ws.Range("J3").WrapText = $true
$ws.Range("M3").Value = "@Get"

# ---------------------------------------------------------------------------
# 3. Row 4 - "Test one- PUT"
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Test one- PUT"
$ws.Range("B4").Value = "REST"
$ws.Range("C4").Value = "ep"
$ws.Range("D4").Value = "Put user"
$ws.Range("E4").Value = "https://reqres.in/api/users/2"
$ws.Range("F4").Value = "PUT"
$ws.Range("G4").Value = "application/json"
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = "name=venom"
$ws.Range("J4").Value = "{`n    ""name"": ""venom"",`n    ""job"": ""zion resident""`n}"
$ws.Range("J4").WrapText = $true
$ws.Range("M4").Value = "@Get"

# ---------------------------------------------------------------------------
# 4. Row 5 - "Test one - PATCH"
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Test one - PATCH"
$ws.Range("B5").Value = "REST"
$ws.Range("C5").Value = "ep"
$ws.Range("D5").Value = "PATCH user"
$ws.Range("E5").Value = "https://reqres.in/api/users/2"
$ws.Range("F5").Value = "PATCH"
$ws.Range("G5").Value = "application/json"
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = "name=venom"
$ws.Range("J5").Value = "{`n    ""name"": ""venom"",`n    ""job"": ""zion resident""`n}"
$ws.Range("J5").WrapText = $true
$ws.Range("M5").Value = "@Get"

# ---------------------------------------------------------------------------
# 5. Row 6 - "Test one - DELETE"
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Test one - DELETE"
$ws.Range("B6").Value = "REST"
$ws.Range("C6").Value = "ep"
$ws.Range("D6").Value = "DELETE user"
$ws.Range("E6").Value = "https://reqres.in/api/users/2"
$ws.Range("F6").Value = "DELETE"
$ws.Range("G6").Value = "application/json"
$ws.Range("H6").Value = 204
$ws.Range("M6").Value = "@Get"

# ---------------------------------------------------------------------------
# 6. Copy the "Tags" cell style (border + number format, used on M2) onto the
#    new M3:M6 cells so they visually match the rest of the @Get tag column.
# ---------------------------------------------------------------------------
$ws.Range("M2").Copy() | Out-Null
$ws.Range("M3:M6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 7. Row heights for the wrapped JSON rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 57.6
$ws.Rows.Item(4).RowHeight = 57.6
$ws.Rows.Item(5).RowHeight = 57.6

# ---------------------------------------------------------------------------
# 8. Column widths - widen the new RequestContent column, and re-apply the
#    (slightly adjusted) widths of the other data columns.
# ---------------------------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 24.5546875

# ---------------------------------------------------------------------------
# 9. View state - selection moved to K5, scrolled so column F is leftmost.
# ---------------------------------------------------------------------------
$ws.Range("K5").Select()

Write-Host "Regres.xlsx updated with POST/PUT/PATCH/DELETE API test rows"
